$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-style rows 8 and 9 (MiniUSB / SIP connector rows) ---------------
# They currently use the now-redundant xf indices 9/10 (quotePrefix/no-quotePrefix
# variants of fillId=0,borderId=1). Those are duplicates of xf 2/1 respectively
# (already used throughout the sheet, e.g. rows 23/24), so we re-point the
# formatting of rows 8 & 9 at the existing xf 2/1 pair by copying the cell
# *formats only* from row 23, which already has exactly that style pattern
# (2 2 2 2 2 1 1 1 2 1) - leaving the duplicate xf entries unused so they can
# be dropped from cellXfs on save.
$ws.Range("A23:J23").Copy() | Out-Null
$ws.Range("A8:J8").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:J23").Copy() | Out-Null
$ws.Range("A9:J9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Add the reset MOSFET row (row 25) ----------------------------------
# Full-row style is a uniform xf 2 (border, no fill) across all 10 columns -
# same pattern as row 12's trailing quantity/blank cells. Seed the row format
# from row 18 (xf 2/1 pattern) then override F to xf 2 as well.
$ws.Range("A18:J18").Copy() | Out-Null
$ws.Range("A25:J25").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Copy() | Out-Null
$ws.Range("F25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A25").Value = "2N7000_TO92"
$ws.Range("B25").Value = "MOSFET N-CH 60V 200MA TO-92"
$ws.Range("C25").Value = "Q1, Q2"
$ws.Range("D25").Value = "TO92_DGS"
$ws.Range("E25").Value = "MOSFET N"
$ws.Range("F25").Value = 2

# --- 3. Add the reset resistor row (row 26) --------------------------------
# Style pattern (2 2 2 2 2 1 2 2 2 2) matches row 18/19 exactly.
$ws.Range("A18:J18").Copy() | Out-Null
$ws.Range("A26:J26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A26").Value = "4K7"
$ws.Range("B26").Value = "RES0805 4K7 5%"
$ws.Range("C26").Value = "R13, R14"
$ws.Range("D26").Value = "0805"
$ws.Range("E26").Value = "RES0805_4K7_5%"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = "Stackpole Electronics Inc"
$ws.Range("H26").Value = "RMCF0805JT4K70"
$ws.Range("I26").Value = "Digi-Key"
$ws.Range("J26").Value = "RMCF0805JT4K70TR-ND"

# --- 4. Selection moves to D27, like the recorded edit ---------------------
$ws.Range("D27").Select() | Out-Null
